$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "29.119.20"
$ws.Range("E2").Value2 = "  -2.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.897.09"
$ws.Range("E3").Value2 = "  -2.79%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.002"
$ws.Range("E4").Value2 = "  +0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "331.36"
$ws.Range("E5").Value2 = "  -3.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "1.001"
$ws.Range("E6").Value2 = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4605"
$ws.Range("E7").Value2 = "  -3.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.4129"
$ws.Range("E8").Value2 = "  -0.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "47.80"
$ws.Range("E9").Value2 = "  -1.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.07994"
$ws.Range("E10").Value2 = "  -3.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "1.008"
$ws.Range("E11").Value2 = "  -3.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "22.26"
$ws.Range("E12").Value2 = "  -1.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "1.904.40"
$ws.Range("E13").Value2 = "  -2.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "5.939"
$ws.Range("E14").Value2 = "  -4.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "7.106"
$ws.Range("E15").Value2 = "  -4.41%  "

$ws.Range("B16").Value2 = "BinanceUSD"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "1.003"
$ws.Range("E16").Value2 = "  +0.20%  "

$ws.Range("B17").Value2 = "Litecoin"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "89.14"
$ws.Range("E17").Value2 = "  -3.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.00001030"
$ws.Range("E18").Value2 = "  -3.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.06564"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "17.62"
$ws.Range("E20").Value2 = "  -2.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "1.003"
$ws.Range("E21").Value2 = "  +0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "29.112.50"
$ws.Range("E22").Value2 = "  -2.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.469"
$ws.Range("E23").Value2 = "  -2.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "11.43"
$ws.Range("E24").Value2 = "  +1.25%  "

$ws.Range("E25").Value2 = "  -3.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "2.135.87"
$ws.Range("E26").Value2 = "  -1.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "156.62"
$ws.Range("E27").Value2 = "  -2.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "19.71"
$ws.Range("E28").Value2 = "  -2.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.118"
$ws.Range("E29").Value2 = "  -3.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "5.635"
$ws.Range("E30").Value2 = "  -1.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "117.34"
$ws.Range("E31").Value2 = "  -4.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.039"
$ws.Range("E32").Value2 = "  +1.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.09378"
$ws.Range("E33").Value2 = "  -2.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.415"
$ws.Range("E34").Value2 = "  -4.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "3.523"
$ws.Range("E35").Value2 = "  -4.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "5.350"
$ws.Range("E36").Value2 = "  -3.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.06073"
$ws.Range("E37").Value2 = "  -3.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.02237"
$ws.Range("E38").Value2 = "  -4.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "8.416"
$ws.Range("E39").Value2 = "  -1.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.173"
$ws.Range("E40").Value2 = "  -2.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.5832"
$ws.Range("E41").Value2 = "  -4.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "1.002"
$ws.Range("E42").Value2 = "  +0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.1825"
$ws.Range("E43").Value2 = "  -3.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "10.14"
$ws.Range("E44").Value2 = "  -5.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "1.266"
$ws.Range("E45").Value2 = "  -1.84%  "

$ws.Range("B46").Value2 = "Cronos"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.07526"
$ws.Range("E46").Value2 = "  +1.73%  "

$ws.Range("B47").Value2 = "RenderToken"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "2.297"
$ws.Range("E47").Value2 = "  -4.05%  "

$ws.Range("B48").Value2 = "EnergySwap"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "12.10"
$ws.Range("E48").Value2 = "  -3.60%  "

$ws.Range("B49").Value2 = "Decentraland"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.5511"
$ws.Range("E49").Value2 = "  -3.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.922"
$ws.Range("E50").Value2 = "  -4.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "111.92"
$ws.Range("E51").Value2 = "  -2.10%  "
